$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append three new rows (290-292) of time/cost log data to the sheet,
# mirroring the existing "time" (text) / "cost" (number) column layout.
$ws.Range("A290").Value = "2023-12-13 18:16:10"
$ws.Range("B290").Value = 0.0012

$ws.Range("A291").Value = "2023-12-13 18:17:10"
$ws.Range("B291").Value = 0.003200000000000001

$ws.Range("A292").Value = "2023-12-13 18:17:36"
$ws.Range("B292").Value = 0.0008
